# Apply updated cryptocurrency price/volume figures to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# --- Price column (D) updates ---
$ws.Range("D2").Value = '61.018.52'
$ws.Range("D3").Value = '2.923.34'
Set-TextValue $ws.Range("D5") '587.11'
Set-TextValue $ws.Range("D6") '146.17'
$ws.Range("D9").Value = '2.922.12'
Set-TextValue $ws.Range("D10") '6.86'
Set-TextValue $ws.Range("D14") '33.68'
$ws.Range("D16").Value = '3.407.18'
$ws.Range("D17").Value = '60.976.38'
$ws.Range("D19").Value = '2.923.34'
Set-TextValue $ws.Range("D20") '431.23'
Set-TextValue $ws.Range("D21") '13.62'
Set-TextValue $ws.Range("D22") '0.682'
Set-TextValue $ws.Range("D23") '7.14'
Set-TextValue $ws.Range("D24") '80.55'
Set-TextValue $ws.Range("D25") '10.84'
Set-TextValue $ws.Range("D30") '7.21'
$ws.Range("D35").Value = '0.0₃0874'
Set-TextValue $ws.Range("D36") '1.01'
Set-TextValue $ws.Range("D40") '49.52'
Set-TextValue $ws.Range("D41") '2.01'
Set-TextValue $ws.Range("D42") '8.67'
Set-TextValue $ws.Range("D43") '0.298'
Set-TextValue $ws.Range("D44") '41.48'
Set-TextValue $ws.Range("D45") '378.92'
$ws.Range("D47").Value = '2.703.14'
Set-TextValue $ws.Range("D48") '132.64'
Set-TextValue $ws.Range("D50") '24.97'

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = '  -2.97%  '
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("E6").Value = '  -4.89%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -2.00%  '
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("E11").Value = '  -4.30%  '
$ws.Range("E12").Value = '  -3.43%  '
$ws.Range("E13").Value = '  -3.46%  '
$ws.Range("E14").Value = '  -5.39%  '
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("E18").Value = '  -4.00%  '
$ws.Range("E19").Value = '  -3.69%  '
$ws.Range("E20").Value = '  -4.88%  '
$ws.Range("E21").Value = '  -4.50%  '
$ws.Range("E22").Value = '  -2.16%  '
$ws.Range("E23").Value = '  -4.87%  '
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("E25").Value = '  -3.76%  '
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -4.21%  '
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("E33").Value = '  -3.42%  '
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("E37").Value = '  -4.56%  '
$ws.Range("E38").Value = '  -5.20%  '
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("E41").Value = '  -5.10%  '
$ws.Range("E42").Value = '  -4.57%  '
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("E44").Value = '  -2.67%  '
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("E51").Value = '  -1.87%  '
